$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 87.333336
$ws.Range("I2").Value = 87.333336
$ws.Range("K2").Value = 87.333336
$ws.Range("M2").Value = 25.666664
$ws.Range("H9").Value = 77.84999999999999
$ws.Range("I9").Value = 23.833334
$ws.Range("K9").Value = 23.833334
$ws.Range("M9").Value = 145.166666
$ws.Range("H53").Value = 1132.2273
$ws.Range("J53").Value = 1249.091
$ws.Range("L53").Value = 1249.091
$ws.Range("N53").Value = -2523.091
$ws.Range("H69").Value = 9631.904
$ws.Range("I69").Value = 6000
$ws.Range("J69").Value = 10014.211
$ws.Range("K69").Value = 18000
$ws.Range("L69").Value = 30042.633
$ws.Range("M69").Value = -17126
$ws.Range("N69").Value = -31790.633
$ws.Range("H72").Value = 9631.904
$ws.Range("I72").Value = 6000
$ws.Range("J72").Value = 10014.211
$ws.Range("K72").Value = 54000
$ws.Range("L72").Value = 90127.89899999999
$ws.Range("M72").Value = -49632
$ws.Range("N72").Value = -98863.89899999999
$ws.Range("H98").Value = 1079.4584
$ws.Range("I98").Value = 1092.8723
$ws.Range("K98").Value = 1092.8723
$ws.Range("M98").Value = 405.1277
$ws.Range("H122").Value = 1079.4584
$ws.Range("I122").Value = 1092.8723
$ws.Range("K122").Value = 3278.6169
$ws.Range("M122").Value = -828.6169
$ws.Range("H124").Value = 562000
$ws.Range("J124").Value = 562000
$ws.Range("L124").Value = 562000
$ws.Range("N124").Value = -571820
$ws.Range("H131").Value = 8463588
$ws.Range("I131").Value = 10000874
$ws.Range("K131").Value = 30002622
$ws.Range("M131").Value = -29997582
$ws.Range("H132").Value = 3369.7112
$ws.Range("I132").Value = 3415.95
$ws.Range("J132").Value = 2999.8
$ws.Range("K132").Value = 10247.85
$ws.Range("L132").Value = 8999.400000000001
$ws.Range("M132").Value = -7717.849999999999
$ws.Range("N132").Value = -14059.4
$ws.Range("H137").Value = 49722.555
$ws.Range("I137").Value = 66079
$ws.Range("J137").Value = 7195.8
$ws.Range("K137").Value = 198237
$ws.Range("L137").Value = 21587.4
$ws.Range("M137").Value = -195687
$ws.Range("N137").Value = -26687.4
$ws.Range("H141").Value = 8500
$ws.Range("I141").Value = 6333.3335
$ws.Range("K141").Value = 19000.0005
$ws.Range("M141").Value = -13820.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8231.799999999999
$ws.Range("I32").Value = 2824.75
$ws.Range("K32").Value = 2824.75
$ws.Range("M32").Value = -2537.75
$ws.Range("H55").Value = 15024
$ws.Range("I55").Value = 15024
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 15024
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -14709
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 2357.2693
$ws.Range("I61").Value = 2156.1667
$ws.Range("K61").Value = 2156.1667
$ws.Range("M61").Value = -1944.1667
$ws.Range("H110").Value = 1303.3334
$ws.Range("I110").Value = 539.375
$ws.Range("J110").Value = 2831.25
$ws.Range("K110").Value = 539.375
$ws.Range("L110").Value = 2831.25
$ws.Range("M110").Value = 1505.625
$ws.Range("N110").Value = -6921.25
$ws.Range("H132").Value = 3024.8718
$ws.Range("I132").Value = 2678.3215
$ws.Range("J132").Value = 3907
$ws.Range("K132").Value = 8034.9645
$ws.Range("L132").Value = 11721
$ws.Range("M132").Value = -5504.9645
$ws.Range("N132").Value = -16781
$ws.Range("H136").Value = 2357.2693
$ws.Range("I136").Value = 2156.1667
$ws.Range("K136").Value = 6468.500100000001
$ws.Range("M136").Value = -3918.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13336455
$ws.Range("I134").Value = 2669
$ws.Range("K134").Value = 8007
$ws.Range("M134").Value = -5472

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3481.9707
$ws.Range("I31").Value = 1867.3684
$ws.Range("J31").Value = 5527.1333
$ws.Range("K31").Value = 1867.3684
$ws.Range("L31").Value = 5527.1333
$ws.Range("M31").Value = -1572.3684
$ws.Range("N31").Value = -6117.1333
$ws.Range("H34").Value = 3481.9707
$ws.Range("I34").Value = 1867.3684
$ws.Range("J34").Value = 5527.1333
$ws.Range("K34").Value = 1867.3684
$ws.Range("L34").Value = 5527.1333
$ws.Range("M34").Value = -1665.3684
$ws.Range("N34").Value = -5931.1333
$ws.Range("H58").Value = 2696.276
$ws.Range("J58").Value = 3687.7778
$ws.Range("L58").Value = 3687.7778
$ws.Range("N58").Value = -4093.7778
$ws.Range("H94").Value = 2509.5264
$ws.Range("I94").Value = 1603.6666
$ws.Range("J94").Value = 2679.375
$ws.Range("K94").Value = 1603.6666
$ws.Range("L94").Value = 2679.375
$ws.Range("M94").Value = -1152.6666
$ws.Range("N94").Value = -3581.375
$ws.Range("H134").Value = 3153.3462
$ws.Range("I134").Value = 2210.2632
$ws.Range("K134").Value = 6630.7896
$ws.Range("M134").Value = -4095.7896
$ws.Range("H136").Value = 2696.276
$ws.Range("J136").Value = 3687.7778
$ws.Range("L136").Value = 11063.3334
$ws.Range("N136").Value = -16163.3334
$ws.Range("H141").Value = 980721.4399999999
$ws.Range("J141").Value = 1162233.2
$ws.Range("L141").Value = 1162233.2
$ws.Range("N141").Value = -1172593.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 598.5
$ws.Range("I34").Value = 598.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1795.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1711.5
$ws.Range("N34").ClearContents()
$ws.Range("H68").Value = 1496.5
$ws.Range("I68").Value = 1495.3334
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 4486.0002
$ws.Range("L68").Value = 4500
$ws.Range("M68").Value = -3675.0002
$ws.Range("N68").Value = -6122
$ws.Range("H71").Value = 1496.5
$ws.Range("I71").Value = 1495.3334
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 13458.0006
$ws.Range("L71").Value = 13500
$ws.Range("M71").Value = -9402.000599999999
$ws.Range("N71").Value = -21612
$ws.Range("H98").Value = 946.7273
$ws.Range("I98").Value = 947.5
$ws.Range("J98").Value = 946.55554
$ws.Range("K98").Value = 2842.5
$ws.Range("L98").Value = 2839.66662
$ws.Range("M98").Value = -1344.5
$ws.Range("N98").Value = -5835.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3331.6667
$ws.Range("J80").Value = 3499.5
$ws.Range("L80").Value = 3499.5
$ws.Range("N80").Value = -5495.5
$ws.Range("H83").Value = 3331.6667
$ws.Range("J83").Value = 3499.5
$ws.Range("L83").Value = 17497.5
$ws.Range("N83").Value = -27481.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 28002.334
$ws.Range("I25").Value = 28002.334
$ws.Range("K25").Value = 28002.334
$ws.Range("M25").Value = -27772.334
$ws.Range("H40").Value = 7280.391
$ws.Range("I40").Value = 7021.4287
$ws.Range("K40").Value = 7021.4287
$ws.Range("M40").Value = -6885.4287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 25006.5
